$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold + border) from an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = 93   # AD
    $ws.Cells.Item($r, 31).Value = 69   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
